# Natmi following Dr Hou advice
# Update the LR-pair stats (columns E:T, rows 2-11) of Anxa2-Robo4 sheet
# to reflect recalculated values (ligand-expressing-cells count changed
# from 1 to 3, etc.), as produced by the updated NATMI pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,125.4603323333333,376.380997,0.1427077914074319,0.1427077914074319,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,6189.140153305327,55702.26137974794,0.1425642353023634,0.1425642353023635),
    @(3,1,125.4603323333333,376.380997,0.1427077914074319,0.1427077914074319,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,6.232200188547555,56.08980169692799,0.0001435561050684928,0.0001435561050684929),
    @(3,1,207.849991,623.549973,0.2364238369850373,0.2364238369850374,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,10253.54150780028,92281.87357020253,0.2361860077477673,0.2361860077477673),
    @(3,1,207.849991,623.549973,0.2364238369850373,0.2364238369850374,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,10.32487901959467,92.923911176352,0.0002378292372700311,0.0002378292372700312),
    @(3,1,240.299296,720.897888,0.2733340584321854,0.2733340584321856,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,11854.31278575896,106688.8150718307,0.2730590995639688,0.2730590995639689),
    @(3,1,240.299296,720.897888,0.2733340584321854,0.2733340584321856,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,11.93678742903467,107.431086861312,0.0002749588682166719,0.000274958868216672),
    @(3,1,190.0452853333333,570.135856,0.2161714578059467,0.2161714578059467,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,9375.209554505494,84376.88599054945,0.2159540013363065,0.2159540013363066),
    @(3,1,190.0452853333333,570.135856,0.2161714578059467,0.2161714578059467,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,9.44043620050489,84.96392580454399,0.0002174564696401267,0.0002174564696401267),
    @(3,1,115.486529,346.459587,0.1313628553693985,0.1313628553693985,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,5697.117969532559,51274.06172579304,0.1312307116393224,0.1312307116393224),
    @(3,1,115.486529,346.459587,0.1313628553693985,0.1313628553693985,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,5.736754832565333,51.630793493088,0.0001321437300761458,0.0001321437300761458)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        # columns E(5) .. T(20)
        $colNum = $j + 5
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
